# İş Takip Güncellemesi - 21.11.2025 08:53:21
#
# 1) "İş Takip Listesi" sheet: columns J (UÇUŞ/İŞE BAŞLAMA benzeri tarih) and
#    K shift back by one day for every populated row (2-122).
# 2) Same sheet: append a new progress-note line (dated 20.11.2025) onto the
#    existing multi-line text in column M for a handful of rows.
# 3) "Güncelleme" sheet: columns I, J, N and P shift back by one day for
#    every populated row (2-29).
#
# Helper that stores a value as literal text (never let Excel's
# autodetect turn a yyyy-MM-dd-looking string into a real date serial),
# while keeping the cell's style index untouched (0 / "Normal").
function Set-TextValue($range, [string]$val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

# Subtract one calendar day from a "yyyy-MM-dd" string; empty input passes
# straight through untouched.
function Get-ShiftedDate([string]$s) {
    if ([string]::IsNullOrEmpty($s)) {
        return $s
    }
    $d = [DateTime]::ParseExact($s, "yyyy-MM-dd", $null)
    return $d.AddDays(-1).ToString("yyyy-MM-dd")
}

function Shift-CellByOneDay($ws, [int]$row, [int]$col) {
    $cell = $ws.Cells.Item($row, $col)
    $cur = $cell.Value()
    if ([string]::IsNullOrEmpty($cur)) {
        return
    }
    $newval = Get-ShiftedDate $cur
    Set-TextValue $cell $newval
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: İş Takip Listesi
# ---------------------------------------------------------------------------
$wsIs = $wb.Worksheets.Item("İş Takip Listesi")

for ($r = 2; $r -le 122; $r++) {
    Shift-CellByOneDay $wsIs $r 10   # column J
    Shift-CellByOneDay $wsIs $r 11   # column K
}

# Append follow-up notes (column M) for the rows that got a new update line.
$notesToAppend = @{
    102 = "20.11.2025 orman mükerrerliği ile ilgili beyanname düzenlendi tescili bekleniyor"
    106 = "20.11.2025 arazi kontrolü yapılıyor 24.11.2025 de firmaya teslim edilecek"
    111 = "20.11.2025 ormancı kontrolünde"
    117 = "20.11.2025 krokilere 22.11.2025 de başlanarak 24.11.2025 kontrol için verilecek"
    118 = "20.11.2025 firmaya 24.11.2025 de teslim edilecek"
}

foreach ($row in $notesToAppend.Keys) {
    $cell = $wsIs.Cells.Item($row, 13)   # column M
    $existing = $cell.Value()
    $cell.Value = $existing + "`n" + $notesToAppend[$row]
}

# ---------------------------------------------------------------------------
# Sheet 2: Güncelleme
# ---------------------------------------------------------------------------
$wsGuncelleme = $wb.Worksheets.Item("Güncelleme")

for ($r = 2; $r -le 29; $r++) {
    Shift-CellByOneDay $wsGuncelleme $r 9    # column I
    Shift-CellByOneDay $wsGuncelleme $r 10   # column J
    Shift-CellByOneDay $wsGuncelleme $r 14   # column N
    Shift-CellByOneDay $wsGuncelleme $r 16   # column P
}
